$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Cells.Item(2, 5).Value = "'1627"
$ws.Cells.Item(3, 5).Value = "'1597"
$ws.Cells.Item(4, 5).Value = "'1565"
$ws.Cells.Item(5, 5).Value = "'1492"
$ws.Cells.Item(6, 5).Value = "'1505"
$ws.Cells.Item(7, 5).Value = "'1474"
$ws.Cells.Item(8, 5).Value = "'1444"
$ws.Cells.Item(9, 5).Value = "'1414"
$ws.Cells.Item(10, 5).Value = "'1435"
$ws.Cells.Item(11, 5).Value = "'1457"
$ws.Cells.Item(12, 5).Value = "'1468"
$ws.Cells.Item(13, 5).Value = "'1482"
$ws.Cells.Item(14, 5).Value = "'1400"
$ws.Cells.Item(15, 5).Value = "'1430"
$ws.Cells.Item(16, 5).Value = "'1492"
$ws.Cells.Item(17, 5).Value = "'1532"
$ws.Cells.Item(18, 5).Value = "'1549"
$ws.Cells.Item(19, 5).Value = "'1525"
$ws.Cells.Item(20, 5).Value = "'1545"
$ws.Cells.Item(21, 5).Value = "'1551"
$ws.Cells.Item(22, 5).Value = "'1543"
$ws.Cells.Item(23, 5).Value = "'1535"
$ws.Cells.Item(24, 5).Value = "'1640"
$ws.Cells.Item(25, 5).Value = "'1666"
$ws.Cells.Item(26, 5).Value = "'1733"
$ws.Cells.Item(27, 5).Value = "'1570"
$ws.Cells.Item(28, 5).Value = "'1548"
$ws.Cells.Item(29, 5).Value = "'1570"
$ws.Cells.Item(30, 5).Value = "'1519"
$ws.Cells.Item(31, 5).Value = "'1594"
$ws.Cells.Item(32, 5).Value = "'1691"
$ws.Cells.Item(33, 5).Value = "'1733"
$ws.Cells.Item(34, 5).Value = "'1862"
$ws.Cells.Item(35, 5).Value = "'1728"
$ws.Cells.Item(36, 5).Value = "'1809"
$ws.Cells.Item(37, 5).Value = "'1886"
$ws.Cells.Item(38, 5).Value = "'1867"
$ws.Cells.Item(39, 5).Value = "'1782"
$ws.Cells.Item(40, 5).Value = "'1780"
$ws.Cells.Item(41, 5).Value = "'1715"
$ws.Cells.Item(42, 5).Value = "'1811"
$ws.Cells.Item(43, 5).Value = "'1812.42146924239"
$ws.Cells.Item(44, 5).Value = "'1790.53006580995"
$ws.Cells.Item(45, 5).Value = "'1791.963217413"
$ws.Cells.Item(46, 5).Value = "'1731.80341707125"
$ws.Cells.Item(47, 5).Value = "'1835.57459064437"
$ws.Cells.Item(48, 5).Value = "'1840.48584407706"
$ws.Cells.Item(49, 5).Value = "'1867.60589331984"
$ws.Cells.Item(50, 5).Value = "'1862.7557348701"
$ws.Cells.Item(51, 5).Value = "'1881.69717862755"
$ws.Cells.Item(52, 5).Value = "'1908.93662092612"
$ws.Cells.Item(53, 5).Value = "'1928.09663053527"
$ws.Cells.Item(54, 5).Value = "'1934.39674920161"
$ws.Cells.Item(55, 5).Value = "'1919.61510267392"
$ws.Cells.Item(56, 5).Value = "'1923.71929539071"
$ws.Cells.Item(57, 5).Value = "'1875.91932393658"
$ws.Cells.Item(58, 5).Value = "'1870.53250067736"
$ws.Cells.Item(59, 5).Value = "'1904.54958156751"
$ws.Cells.Item(60, 5).Value = "'1920.30605578009"
$ws.Cells.Item(61, 1).Value = 204
$ws.Cells.Item(61, 2).Value = "Benin"
$ws.Cells.Item(61, 3).Value = "GDP per Capita"
$ws.Cells.Item(61, 4).Value = 2009
$ws.Cells.Item(61, 5).Value = "'1888.77577922789"
$ws.Cells.Item(62, 1).Value = 204
$ws.Cells.Item(62, 2).Value = "Benin"
$ws.Cells.Item(62, 3).Value = "GDP per Capita"
$ws.Cells.Item(62, 4).Value = 2010
$ws.Cells.Item(62, 5).Value = "'1854.66557072809"
$ws.Cells.Item(63, 1).Value = 204
$ws.Cells.Item(63, 2).Value = "Benin"
$ws.Cells.Item(63, 3).Value = "GDP per Capita"
$ws.Cells.Item(63, 4).Value = 2011
$ws.Cells.Item(63, 5).Value = "'1836"
$ws.Cells.Item(64, 1).Value = 204
$ws.Cells.Item(64, 2).Value = "Benin"
$ws.Cells.Item(64, 3).Value = "GDP per Capita"
$ws.Cells.Item(64, 4).Value = 2012
$ws.Cells.Item(64, 5).Value = "'1867"
$ws.Cells.Item(65, 1).Value = 204
$ws.Cells.Item(65, 2).Value = "Benin"
$ws.Cells.Item(65, 3).Value = "GDP per Capita"
$ws.Cells.Item(65, 4).Value = 2013
$ws.Cells.Item(65, 5).Value = "'1939"
$ws.Cells.Item(66, 1).Value = 204
$ws.Cells.Item(66, 2).Value = "Benin"
$ws.Cells.Item(66, 3).Value = "GDP per Capita"
$ws.Cells.Item(66, 4).Value = 2014
$ws.Cells.Item(66, 5).Value = "'2008"
$ws.Cells.Item(67, 1).Value = 204
$ws.Cells.Item(67, 2).Value = "Benin"
$ws.Cells.Item(67, 3).Value = "GDP per Capita"
$ws.Cells.Item(67, 4).Value = 2015
$ws.Cells.Item(67, 5).Value = "'2055"
$ws.Cells.Item(68, 1).Value = 204
$ws.Cells.Item(68, 2).Value = "Benin"
$ws.Cells.Item(68, 3).Value = "GDP per Capita"
$ws.Cells.Item(68, 4).Value = 2016
$ws.Cells.Item(68, 5).Value = "'2080"
